$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "W-000001"
$ws.Range("C2").Value = "B-000024"
$ws.Range("D2").Value = 13
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "COL"
$ws.Range("H2").Value = "Sohel1"
$ws.Range("I2").Value = "Hajipara"
$ws.Range("J2").Value = "Hajipara 2"

$ws.Range("A1").Copy()
$ws.Range("B2:J2").PasteSpecial(-4122)
